# semana 48 de 2024
# Update the "Esperado" (C), "Observado" (D) and "valor p" (E) columns
# for the epidemiological events table with the latest week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1
$ws.Range("E2").Value = 0

$ws.Range("C3").Value = 3
$ws.Range("D3").Value = 2
$ws.Range("E3").Value = 0.22

$ws.Range("D4").Value = 1
$ws.Range("E4").Value = 0

$ws.Range("C5").Value = 8
$ws.Range("D5").Value = 5
$ws.Range("E5").Value = 0.09

$ws.Range("C6").Value = 2
$ws.Range("D6").Value = 30

$ws.Range("C7").Value = 1
$ws.Range("D7").Value = 6
$ws.Range("E7").Value = 0

$ws.Range("D9").Value = 0
$ws.Range("E9").Value = 1

$ws.Range("C11").Value = 43
$ws.Range("D11").Value = 37

$ws.Range("C12").Value = 0
$ws.Range("E12").Value = 1

$ws.Range("D13").Value = 1

$ws.Range("C14").Value = 3
$ws.Range("D14").Value = 3
$ws.Range("E14").Value = 0.22

$ws.Range("C15").Value = 6
$ws.Range("E15").Value = 0

$ws.Range("C16").Value = 0
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 0

$ws.Range("C17").Value = 0
$ws.Range("E17").Value = 1

$ws.Range("C18").Value = 4
$ws.Range("E18").Value = 0.02

$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 0.14

$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 0.27

$ws.Range("C21").Value = 8
$ws.Range("D21").Value = 0
$ws.Range("E21").Value = 0

$ws.Range("C24").Value = 2
$ws.Range("D24").Value = 3
$ws.Range("E24").Value = 0.18

$ws.Range("C25").Value = 1
$ws.Range("E25").Value = 0.37

$ws.Range("C27").Value = 6
$ws.Range("D27").Value = 6
$ws.Range("E27").Value = 0.16

$ws.Range("D28").Value = 1
$ws.Range("E28").Value = 0.37

$ws.Range("D29").Value = 0
$ws.Range("E29").Value = 1

$ws.Range("C30").Value = 1
$ws.Range("E30").Value = 0.37

$ws.Range("C34").Value = 2
$ws.Range("D34").Value = 1
$ws.Range("E34").Value = 0.27

$ws.Range("C35").Value = 6
$ws.Range("D35").Value = 4
$ws.Range("E35").Value = 0.13

$ws.Range("C36").Value = 7
$ws.Range("D36").Value = 4
$ws.Range("E36").Value = 0.09

$ws.Range("C37").Value = 9
$ws.Range("D37").Value = 6
$ws.Range("E37").Value = 0.09
